$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.102.12"
Set-TextValue "E2" "  -2.53%  "
Set-TextValue "D3" "1.866.70"
Set-TextValue "E3" "  -2.04%  "
Set-TextValue "D4" "0.9999"
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "306.22"
Set-TextValue "E5" "  -2.04%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "D7" "0.5167"
Set-TextValue "E7" "  -1.33%  "
Set-TextValue "D8" "0.3767"
Set-TextValue "D9" "0.07164"
Set-TextValue "E9" "  -0.95%  "
Set-TextValue "E10" "  -1.89%  "
Set-TextValue "E11" "  -2.74%  "
Set-TextValue "D12" "0.07602"
Set-TextValue "E12" "  -0.63%  "
Set-TextValue "D13" "1.849.22"
Set-TextValue "E13" "  -2.92%  "
Set-TextValue "E14" "  -2.52%  "
Set-TextValue "D15" "89.78"
Set-TextValue "E15" "  -2.51%  "
Set-TextValue "D16" "1.000"
Set-TextValue "E16" "  -0.05%  "
Set-TextValue "D17" "0.000008487"
Set-TextValue "E18" "  -3.19%  "
Set-TextValue "E19" "  +0.00%  "
Set-TextValue "D20" "27.127.97"
Set-TextValue "E20" "  -2.54%  "
Set-TextValue "D21" "5.033"
Set-TextValue "E21" "  -2.31%  "
Set-TextValue "D22" "2.116.21"
Set-TextValue "E22" "  -1.54%  "
Set-TextValue "E23" "  -3.24%  "
Set-TextValue "D24" "6.468"
Set-TextValue "E24" "  -2.47%  "
Set-TextValue "E25" "  -1.73%  "
Set-TextValue "D26" "147.55"
Set-TextValue "E26" "  -3.95%  "
Set-TextValue "D27" "17.97"
Set-TextValue "E27" "  -1.97%  "
Set-TextValue "D28" "2.099"
Set-TextValue "E28" "  -3.13%  "
Set-TextValue "D29" "112.94"
Set-TextValue "E29" "  -1.55%  "
Set-TextValue "D30" "4.666"
Set-TextValue "E30" "  -4.03%  "
Set-TextValue "D31" "4.673"
Set-TextValue "E31" "  -3.68%  "
Set-TextValue "D32" "0.09148"
Set-TextValue "E32" "  +1.36%  "
Set-TextValue "D33" "0.05118"
Set-TextValue "E33" "  -3.00%  "
Set-TextValue "D34" "3.071"
Set-TextValue "E34" "  -3.31%  "
Set-TextValue "D35" "1.159"
Set-TextValue "E35" "  -6.18%  "
Set-TextValue "D36" "0.7276"
Set-TextValue "E36" "  -6.75%  "
Set-TextValue "E37" "  -2.83%  "
Set-TextValue "D38" "3.075"
Set-TextValue "E38" "  +0.17%  "
Set-TextValue "D39" "2.503"
Set-TextValue "E39" "  -4.22%  "
Set-TextValue "D40" "1.075"
Set-TextValue "E40" "  -1.61%  "
Set-TextValue "D41" "0.5326"
Set-TextValue "E41" "  -4.76%  "
Set-TextValue "D42" "6.485"
Set-TextValue "E42" "  -3.47%  "
Set-TextValue "D43" "116.25"
Set-TextValue "D44" "8.291"
Set-TextValue "E44" "  -3.33%  "
Set-TextValue "D45" "0.1466"
Set-TextValue "E45" "  -3.20%  "
Set-TextValue "D46" "0.4637"
Set-TextValue "E46" "  -3.38%  "
Set-TextValue "D47" "0.9996"
Set-TextValue "E47" "  +0.04%  "
Set-TextValue "D48" "9.969"
Set-TextValue "E48" "  -5.03%  "
Set-TextValue "D49" "1.570"
Set-TextValue "E49" "  -3.26%  "
Set-TextValue "D50" "36.59"
Set-TextValue "E50" "  -1.10%  "
Set-TextValue "D51" "63.62"
Set-TextValue "E51" "  -4.94%  "
